# Move the "detection ceiling" value out of the time_variants scatter and
# into the constants/parameters sheet as a brand new scalar parameter
# (tb_prop_detection_algorithm_ceiling = 0.95).

$wb = $excel.ActiveWorkbook
$wsConstants = $wb.Worksheets.Item("constants")
$wsTimeVariants = $wb.Worksheets.Item("time_variants")

# ---------------------------------------------------------------------
# 1. constants sheet: insert a new row above the old "epi_prop_smearpos"
#    row (row 38) and populate it with the new parameter. Everything
#    below shifts down by one row automatically.
# ---------------------------------------------------------------------
$wsConstants.Rows("38:38").Insert()
$wsConstants.Range("A38").Value2 = "tb_prop_detection_algorithm_ceiling"
$wsConstants.Range("B38").Value2 = 0.95

# ---------------------------------------------------------------------
# 2. time_variants sheet: the "program_perc_detect" row (row 3) had data
#    points at 1930 (col F) and 1965 (col K). Move them to 1950 (col H)
#    and 1960 (col J) respectively.
# ---------------------------------------------------------------------
$wsTimeVariants.Range("H3").Value2 = $wsTimeVariants.Range("F3").Value2
$wsTimeVariants.Range("F3").Value2 = $null
$wsTimeVariants.Range("J3").Value2 = $wsTimeVariants.Range("K3").Value2
$wsTimeVariants.Range("K3").Value2 = $null

# ---------------------------------------------------------------------
# 3. View/selection bookkeeping to mirror the author's saved state:
#    active sheet moves from time_variants to constants, with the
#    selection sitting on the freshly edited row.
# ---------------------------------------------------------------------
$wsTimeVariants.Activate()
$wsTimeVariants.Range("C16").Select()

$wsConstants.Activate()
$wsConstants.Range("A39").Select()
